$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths for new columns E (5) and F (6) ---
# Widths are stored internally in 1/6-character increments by this engine, so we
# pick the ColumnWidth value whose rounded stored width lands on the nearest
# achievable grid point to the target (13.21875 / 10.109375 chars).
$ws.Columns.Item(5).ColumnWidth = 12.26
$ws.Columns.Item(6).ColumnWidth = 9.26

# --- New rows 12 & 13: copy formatting from row 11 (same style pattern: col A = s2, B:L = s3) ---
$ws.Range("A11:L11").Copy()
$ws.Range("A12:L13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 12 values ---
$ws.Range("A12").Value = "TC_11_Validate Login Page Links"
$ws.Range("B12").Value = "admin"
$ws.Range("C12").Value = "admin"
$ws.Range("D12").Value = "Aqua"
$ws.Range("E12").Value = "Senior QA Lead"
$ws.Range("F12").Value = "QA_@123"

# --- Row 13 values ---
$ws.Range("A13").Value = "TC_12_Validate New Vendor Page"
$ws.Range("B13").Value = "admin"
$ws.Range("C13").Value = "admin"
$ws.Range("D13").Value = "nature"
$ws.Range("E13").Value = "Senior QA Lead"
$ws.Range("F13").Value = "QA_@123"

# --- Hyperlinks on F12/F13 (adds Hyperlink font/style + cellXfs entry automatically) ---
$ws.Hyperlinks.Add($ws.Range("F12"), "https://example.com")
$ws.Hyperlinks.Add($ws.Range("F13"), "https://example.com")

# --- Selection moves from L18 to F18 ---
[void]$ws.Range("F18").Select()

Write-Output "done"
